# Weekly refresh of the "Fruta, Mercado Mayorista Lo Valledor de Santiago - Mora"
# data set: three new weekly observations are inserted into the existing table,
# pushing the later rows down.
#
#   - two new rows are inserted right after row 40 (becoming rows 41 & 42)
#   - one new row is inserted after what is (at that point) row 62,
#     i.e. right after the former row 60's data (becoming row 63)
#
# All rows below each insertion point shift down accordingly, ending with the
# table spanning A1:T65 (was A1:T62).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the first two new rows (after row 40) -------------------------
$ws.Rows("41:42").Insert()

$ws.Range("A41").Value = 6
$ws.Range("B41").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C41").Value = "Metropolitana"
$ws.Range("D41").Value = 44567
$ws.Range("E41").Value = 13
$ws.Range("F41").Value = "Fruta"
$ws.Range("G41").Value = 100101
$ws.Range("H41").Value = "Berries"
$ws.Range("I41").Value = 100101008
$ws.Range("J41").Value = "Mora"
$ws.Range("K41").Value = "Sin especificar"
$ws.Range("L41").Value = "Primera"
$ws.Range("M41").Value = 250
$ws.Range("N41").Value = 6000
$ws.Range("O41").Value = 6000
$ws.Range("P41").Value = 6000
$ws.Range("Q41").Value = "$/bandeja 2 kilos"
$ws.Range("R41").Value = "Provincia de Curicó"
$ws.Range("S41").Value = 3000
$ws.Range("T41").Value = 2

$ws.Range("A42").Value = 6
$ws.Range("B42").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C42").Value = "Metropolitana"
$ws.Range("D42").Value = 44567
$ws.Range("E42").Value = 13
$ws.Range("F42").Value = "Fruta"
$ws.Range("G42").Value = 100101
$ws.Range("H42").Value = "Berries"
$ws.Range("I42").Value = 100101008
$ws.Range("J42").Value = "Mora"
$ws.Range("K42").Value = "Sin especificar"
$ws.Range("L42").Value = "Segunda"
$ws.Range("M42").Value = 250
$ws.Range("N42").Value = 4000
$ws.Range("O42").Value = 4000
$ws.Range("P42").Value = 4000
$ws.Range("Q42").Value = "$/bandeja 2 kilos"
$ws.Range("R42").Value = "Provincia de Curicó"
$ws.Range("S42").Value = 2000
$ws.Range("T42").Value = 2

# --- Insert the third new row (after the row that is now row 62) ----------
$ws.Rows("63:63").Insert()

$ws.Range("A63").Value = 6
$ws.Range("B63").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C63").Value = "Metropolitana"
$ws.Range("D63").Value = 44568
$ws.Range("E63").Value = 13
$ws.Range("F63").Value = "Fruta"
$ws.Range("G63").Value = 100101
$ws.Range("H63").Value = "Berries"
$ws.Range("I63").Value = 100101008
$ws.Range("J63").Value = "Mora"
$ws.Range("K63").Value = "Sin especificar"
$ws.Range("L63").Value = "Primera"
$ws.Range("M63").Value = 250
$ws.Range("N63").Value = 6000
$ws.Range("O63").Value = 6000
$ws.Range("P63").Value = 6000
$ws.Range("Q63").Value = "$/bandeja 2 kilos"
$ws.Range("R63").Value = "Provincia de Linares"
$ws.Range("S63").Value = 3000
$ws.Range("T63").Value = 2
